# 29th March 2022 changes: update two tracking-number cells (P2, P3) on Sheet1
# to new shared-string values, replacing the previous ones.
# The leading apostrophe forces Excel to store the digit string as text
# (t="s" shared-string) instead of coercing it to a number; re-applying the
# "Normal" style afterwards clears the transient quote-prefix formatting so
# the cell's style index is left unchanged, matching a plain value edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = "'320018151874"
$ws.Range("P2").Style = "Normal"

$ws.Range("P3").Value = "'320018151885"
$ws.Range("P3").Style = "Normal"
